$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment the date serial values in column F (rows 2-7) by 1 day each,
# preserving their existing number formatting/style.
for ($row = 2; $row -le 7; $row++) {
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value2 = $cell.Value2 + 1
}
